$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Widen the 4th (Prix) column from 2303 dxa (115.15pt) to 2363 dxa (118.15pt).
# Setting the Column's Width updates the tblGrid entry and every row's tcW
# for that column in one shot.
$t.Columns.Item(4).Width = 118.15

# Append a trailing space after "9€" in row 4 ("PCB Nunchuk" row), column 4,
# as its own run (matching the authored diff which adds a separate <w:r>).
$cell = $t.Cell(4, 4)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Collapse(0)
$rng.InsertAfter(" ")
$newRun = $d.Range($rng.Start, $rng.End)
$newRun.Font.Bold = 1
$newRun.Font.Bold = 0
